$d = $word.ActiveDocument

# Avoid Word's "smart quotes" AutoCorrect from turning the straight
# apostrophe in "don't" into a curly one when we Find/Replace below.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# 1. Remove the "Meta description" paragraph that currently sits right after
#    the title heading (paragraph 2): an empty run, a bold "Meta description"
#    run, and a plain run with the description text.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new paragraph - "Play Book of Sheba for Free - Online Slots
#    Review" in bold - right before the final (italic) paragraph at the end
#    of the document. Build it via raw OOXML so the run layout exactly
#    matches an empty leading run followed by the bold text run, with no
#    stray formatting bleeding in from neighboring paragraphs.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Sheba for Free - Online Slots Review</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newParaXml)

# InsertXML needs a trailing paragraph mark to close off the new paragraph
# properly, which leaves one extra empty paragraph behind - remove it.
$extraPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$extraPara.Range.Delete()

# 3. Replace the text of the final (italic) paragraph with the new review
#    blurb, keeping its existing (italic) character formatting intact.
#    Locate the old text with Find (no Replace argument) and then assign
#    Range.Text directly - using Find's own Replace parameter here would
#    let AutoCorrect turn the straight apostrophe in "don't" into a curly
#    one, which the source diff does not want.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.Find.Execute("Create a Feature Image Prompt: Design an eye-catching, cartoon-style feature image to capture the essence of Book of Sheba. The main focus of the image should be a happy Maya warrior with glasses. The warrior should be dressed in traditional clothing and holding a scepter. The background of the image should feature a desert landscape and ancient Egyptian artifacts like pyramids and hieroglyphs. Use bright, vibrant colors to make the image pop and attract players to this exciting online slot game.")
$finalRange.Text = "Read our review of Book of Sheba online slot game. Play for free and find out what we like and don't like about this exciting Ancient Egypt themed slot."

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
